# Append a new data row (row 29) to the "Artfynd" worksheet, reproducing
# the row that was added to the source OOXML (A 42296-2019.xlsx).
#
# Values that are logically text (species names, dates stored as plain
# strings, etc.) are written with a leading apostrophe so that Excel's
# automatic type-detection (which would otherwise turn "2023-09-09" into
# a real date, or "1" into a number) is suppressed and the value is
# stored as text instead - exactly like the rest of the sheet. The
# apostrophe itself is only an entry-mode marker and is not stored in the
# cell value. After assignment the cell style is reset back to "Normal"
# so no extra/residual number formatting is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

function Set-TextCell {
    param($sheet, $r, $c, [string]$text)
    $cell = $sheet.Cells.Item($r, $c)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

function Set-NumberCell {
    param($sheet, $r, $c, $number)
    $cell = $sheet.Cells.Item($r, $c)
    $cell.Value = $number
}

function Set-BoolCell {
    param($sheet, $r, $c, [bool]$value)
    $cell = $sheet.Cells.Item($r, $c)
    $cell.Value = $value
}

# --- Numeric columns -------------------------------------------------
Set-NumberCell $ws $row 1  111982042                  # A  Id
Set-NumberCell $ws $row 2  56543                       # B  Taxonsorteringsordning
Set-NumberCell $ws $row 5  103021                      # E  TaxonId
Set-NumberCell $ws $row 17 306291.7558343319           # Q  Ost
Set-NumberCell $ws $row 18 6525531.449465405           # R  Nord
Set-NumberCell $ws $row 19 10                          # S  Noggrannhet

# --- Text columns ------------------------------------------------------
Set-TextCell $ws $row 3  "Ovaliderad"                              # C  Valideringsstatus
Set-TextCell $ws $row 4  "NT"                                      # D  Rodlistade
Set-TextCell $ws $row 6  "Talltita"                                # F  Artnamn
Set-TextCell $ws $row 7  "Poecile montanus"                        # G  Vetenskapligt namn
Set-TextCell $ws $row 8  "(Conrad von Baldenstein, 1827)"          # H  Auktor
Set-TextCell $ws $row 9  "1"                                       # I  Antal (stored as text)
Set-TextCell $ws $row 11 ""                                        # K  Alder-Stadium (empty)
Set-TextCell $ws $row 13 "lockläte, övriga läten"                  # M  Aktivitet
Set-TextCell $ws $row 16 "Oxögat (Södra Kornsjön, Kynnefjäll), Boh" # P  Lokalnamn
Set-TextCell $ws $row 20 "Västra Götaland"                         # T  Lan
Set-TextCell $ws $row 21 "Tanum"                                   # U  Kommun
Set-TextCell $ws $row 22 "Bohuslän"                                # V  Provins
Set-TextCell $ws $row 23 "Naverstad"                                # W  Forsamling
Set-TextCell $ws $row 25 "2023-09-09"                               # Y  Startdatum
Set-TextCell $ws $row 26 "13:21"                                    # Z  Starttid
Set-TextCell $ws $row 27 "2023-09-09"                               # AA Slutdatum
Set-TextCell $ws $row 28 "13:21"                                    # AB Sluttid
Set-TextCell $ws $row 46 ""                                         # AT Bestamningsar (empty)
Set-TextCell $ws $row 49 "Christer Johansson"                       # AW Rapportor
Set-TextCell $ws $row 50 "Christer Johansson"                       # AX Observatorer
Set-TextCell $ws $row 51 ""                                         # AY Projektnamn (empty)

# --- Boolean columns ---------------------------------------------------
Set-BoolCell $ws $row 30 $false   # AD Ej aterfunnen
Set-BoolCell $ws $row 31 $false   # AE Osaker artbestamning
Set-BoolCell $ws $row 33 $false   # AG Ospontan
